$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = "dateTime" + [char]0x0135
$ws.Range("I4").Value = "prohibited"

$ws.Range("G7").Value = "dateTime" + [char]0x0135
$ws.Range("I7").Value = "prohibited"
